# Applies the OOXML diff:
#  1. Table 1, row "4" (index cell) gains six new paragraphs describing
#     byte/word/longword storage sizes.
#  2. Table 2 ("Student Name" / "Student Number" / "Date" / "Checked")
#     gets the student's name, number and date filled in.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Table 1 - row whose first cell holds "4": append the new paragraphs
#    after the existing "4" paragraph, preserving its bold/underline-none
#    run formatting for every new paragraph (matching the surrounding
#    cell's formatting).
# ---------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$cell = $t1.Cell(5, 1)
$rpr4 = '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:u w:val="none"/></w:rPr>'

$frag1 = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Normal`"/>$rpr4</w:pPr><w:r>$rpr4<w:t>4</w:t></w:r></w:p>" +
         "<w:p $wNs><w:pPr><w:pStyle w:val=`"Normal`"/>$rpr4</w:pPr></w:p>" +
         "<w:p $wNs><w:pPr><w:pStyle w:val=`"Normal`"/>$rpr4</w:pPr><w:r>$rpr4<w:t>Bytes can store 8 bits</w:t></w:r></w:p>" +
         "<w:p $wNs><w:pPr><w:pStyle w:val=`"Normal`"/>$rpr4</w:pPr><w:r>$rpr4<w:t>Word can store 16 bits</w:t></w:r></w:p>" +
         "<w:p $wNs><w:pPr><w:pStyle w:val=`"Normal`"/>$rpr4</w:pPr><w:r>$rpr4<w:t>Longword stores 32 bits</w:t></w:r></w:p>" +
         "<w:p $wNs><w:pPr><w:pStyle w:val=`"Normal`"/>$rpr4</w:pPr></w:p>" +
         "<w:p $wNs><w:pPr><w:pStyle w:val=`"Normal`"/>$rpr4</w:pPr><w:r>$rpr4<w:t>You can overwrite parts of word of long with bytes</w:t></w:r></w:p>"

$cellRange = $cell.Range
$cellRange.Collapse(0)
$cellRange.InsertXML($frag1)

# ---------------------------------------------------------------------
# 2) Table 2 - fill in the student name / number / date. Each target
#    cell already contains one or more empty paragraphs; the first of
#    those paragraphs gets replaced in place (so any following blank
#    paragraph is left untouched) with a bold run carrying the text.
# ---------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$rprBold = '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr>'

function Fill-FirstParagraph($table, $row, $col, $text, $preserveSpace) {
    $cell = $table.Cell($row, $col)
    $p1 = $cell.Range.Paragraphs.Item(1)
    $rng = $p1.Range
    if ($preserveSpace) {
        $frag = "<w:p $wNs><w:pPr>$rprBold</w:pPr><w:r>$rprBold<w:t xml:space=`"preserve`">$text</w:t></w:r></w:p>"
    } else {
        $frag = "<w:p $wNs><w:pPr>$rprBold</w:pPr><w:r>$rprBold<w:t>$text</w:t></w:r></w:p>"
    }
    $rng.InsertXML($frag)
}

Fill-FirstParagraph $t2 1 2 "Karolis Grigaliunas " $true
Fill-FirstParagraph $t2 1 4 "C00287940" $false
Fill-FirstParagraph $t2 2 2 "22/01" $false
